# Insert a new weekly price record as row 166 in the "Bruselas (repollito)"
# data sheet, pushing the existing rows 166-187 down to 167-188.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 166 (shifts rows 166:187 -> 167:188)
$ws.Rows.Item(166).Insert()

# Populate the new row 166 with the new record's data
$ws.Range("A166").Value = 10
$ws.Range("B166").Value = "Vega Modelo de Temuco"
$ws.Range("C166").Value = "La Araucanía"
$ws.Range("D166").Value = "2023-07-17"
$ws.Range("E166").Value = 9
$ws.Range("F166").Value = 100112035
$ws.Range("G166").Value = "Bruselas (repollito)"
$ws.Range("H166").Value = "Sin especificar"
$ws.Range("I166").Value = "Primera"
$ws.Range("J166").Value = 140
$ws.Range("K166").Value = 25000
$ws.Range("L166").Value = 25000
$ws.Range("M166").Value = 25000
$ws.Range("N166").Value = "$/malla 10 kilos"
$ws.Range("O166").Value = "Provincia de Quillota"
$ws.Range("P166").Value = 1667
$ws.Range("Q166").Value = 15
$ws.Range("R166").Value = "Hortaliza"
